$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3681.335999999998
$ws.Range("C2").Value = 1748.5763141228617
$ws.Range("D2").Value = 9.203339999999995
$ws.Range("E2").Value = 7.079492307692304
$ws.Range("B3").Value = 3534.082559999998
$ws.Range("C3").Value = 1678.6332615579472
$ws.Range("D3").Value = 8.835206399999995
$ws.Range("E3").Value = 6.796312615384611
$ws.Range("B4").Value = 284.8497777777778
$ws.Range("C4").Value = 135.29913447329534
$ws.Range("D4").Value = 0.7121244444444444
$ws.Range("E4").Value = 0.5477880341880342
$ws.Range("B5").Value = 7500.268337777773
$ws.Range("C5").Value = 3562.508710154104
$ws.Range("D5").Value = 18.750670844444432
$ws.Range("E5").Value = 14.42359295726495
$ws.Range("B7").Value = 1515.84
$ws.Range("D7").Value = 3.7895999999999996
$ws.Range("E7").Value = 2.915076923076923
$ws.Range("B8").Value = 568.4399999999999
$ws.Range("D8").Value = 1.4210999999999998
$ws.Range("E8").Value = 1.0931538461538461
$ws.Range("B9").Value = 2084.2799999999997
$ws.Range("D9").Value = 5.210699999999999
$ws.Range("E9").Value = 4.00823076923077
$ws.Range("B11").Value = 1313.497927396606
$ws.Range("C11").Value = 623.8907191560827
$ws.Range("D11").Value = 3.283744818491515
$ws.Range("E11").Value = 2.5259575526857807
$ws.Range("C13").Value = 241.49496516782764
$ws.Range("D13").Value = 1.2710684999999995
$ws.Range("E13").Value = 0.9777449999999996
$ws.Range("C14").Value = 200.8779564285619
$ws.Range("D14").Value = 1.0572876440023307
$ws.Range("E14").Value = 0.8132981876941004
$ws.Range("C15").Value = 629.8290056998098
$ws.Range("D15").Value = 3.314999999999999
$ws.Range("E15").Value = 2.549999999999999
$ws.Range("C16").Value = 1072.201927296199
$ws.Range("D16").Value = 5.643356144002328
$ws.Range("E16").Value = 4.3410431876940985
$ws.Range("B18").Value = 1264.6811855601686
$ws.Range("C18").Value = 600.7035396897571
$ws.Range("D18").Value = 3.1617029639004217
$ws.Range("E18").Value = 2.4320792030003244
$ws.Range("B19").Value = 1739.3130881928244
$ws.Range("C19").Value = 826.1461786856354
$ws.Range("D19").Value = 4.348282720482061
$ws.Range("E19").Value = 3.3448328619092775
$ws.Range("B20").Value = 3115.854584453671
$ws.Range("C20").Value = 1479.9815948956639
$ws.Range("D20").Value = 7.789636461134178
$ws.Range("E20").Value = 5.992028047026291
$ws.Range("B23").Value = 16271.24330722898
$ws.Range("C23").Value = 7728.58295150205
$ws.Range("D23").Value = 40.67810826807245
$ws.Range("E23").Value = 31.290852513901886
$ws.Range("B25").Value = 8770.974969451208
$ws.Range("C25").Value = 4166.074241347946
$ws.Range("D25").Value = 21.927437423628020
$ws.Range("E25").Value = 16.867259556636938
